$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = 5.7948
$ws.Range("C2").Value = 16.33285
$ws.Range("D2").Value = 34.7531
$ws.Range("E2").Value = 6.4706
$ws.Range("F2").Value = 2702.48335
$ws.Range("G2").Value = 2483.29295
$ws.Range("H2").Value = 219.1905
$ws.Range("I2").Value = 219.1905
$ws.Range("K2").Value = 2745.781
$ws.Range("L2").Value = 2483.296
$ws.Range("M2").Value = 262.485
$ws.Range("N2").Value = 262.485
$ws.Range("B3").Value = 5.831
$ws.Range("C3").Value = 16.482
$ws.Range("D3").Value = 34.876
$ws.Range("E3").Value = 5.974
$ws.Range("F3").Value = 2719.272
$ws.Range("G3").Value = 2490.45
$ws.Range("H3").Value = 228.822
$ws.Range("I3").Value = 228.822
$ws.Range("K3").Value = 2744.4092
$ws.Range("L3").Value = 2490.442
$ws.Range("M3").Value = 253.9672
$ws.Range("N3").Value = 253.9672
$ws.Range("B4").Value = 6.158
$ws.Range("C4").Value = 17.154
$ws.Range("D4").Value = 38.034
$ws.Range("E4").Value = 7.722
$ws.Range("F4").Value = 2919.128
$ws.Range("G4").Value = 2591.081
$ws.Range("H4").Value = 328.047
$ws.Range("I4").Value = 328.047
$ws.Range("K4").Value = 2763.7398
$ws.Range("L4").Value = 2591.081
$ws.Range("M4").Value = 172.6588
$ws.Range("N4").Value = 172.6588

$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = 5.78185
$ws.Range("C2").Value = 16.1637
$ws.Range("D2").Value = 35.1861
$ws.Range("E2").Value = 6.3483
$ws.Range("F2").Value = 2690.3983
$ws.Range("G2").Value = 2477.589
$ws.Range("H2").Value = 212.8093
$ws.Range("I2").Value = 212.8093
$ws.Range("K2").Value = 2793.0096
$ws.Range("L2").Value = 2477.597
$ws.Range("M2").Value = 315.4126
$ws.Range("N2").Value = 315.4126
$ws.Range("B3").Value = 5.81
$ws.Range("C3").Value = 16.468
$ws.Range("D3").Value = 34.432
$ws.Range("E3").Value = 5.715
$ws.Range("F3").Value = 2707.623
$ws.Range("G3").Value = 2482.172
$ws.Range("H3").Value = 225.451
$ws.Range("I3").Value = 225.451
$ws.Range("K3").Value = 2792.2678
$ws.Range("L3").Value = 2482.193
$ws.Range("M3").Value = 310.0748
$ws.Range("N3").Value = 310.0748
$ws.Range("B4").Value = 5.983
$ws.Range("C4").Value = 15.185
$ws.Range("D4").Value = 41.485
$ws.Range("E4").Value = 7.372
$ws.Range("F4").Value = 2864.22
$ws.Range("G4").Value = 2518.066
$ws.Range("H4").Value = 346.154
$ws.Range("I4").Value = 346.154
$ws.Range("K4").Value = 2791.2062
$ws.Range("L4").Value = 2518.066
$ws.Range("M4").Value = 273.1402
$ws.Range("N4").Value = 273.1402

$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = 5.79305
$ws.Range("C2").Value = 16.0961
$ws.Range("D2").Value = 36.1616
$ws.Range("E2").Value = 6.843549999999999
$ws.Range("F2").Value = 2709.40465
$ws.Range("G2").Value = 2487.3154
$ws.Range("H2").Value = 222.08925
$ws.Range("I2").Value = 222.08925
$ws.Range("K2").Value = 2707.3304
$ws.Range("L2").Value = 2487.323
$ws.Range("M2").Value = 220.0074
$ws.Range("N2").Value = 220.0074
$ws.Range("B3").Value = 5.8
$ws.Range("C3").Value = 16.249
$ws.Range("D3").Value = 36.172
$ws.Range("E3").Value = 6.657
$ws.Range("F3").Value = 2728.415
$ws.Range("G3").Value = 2492.653
$ws.Range("H3").Value = 235.762
$ws.Range("I3").Value = 235.762
$ws.Range("K3").Value = 2707.8842
$ws.Range("L3").Value = 2492.709
$ws.Range("M3").Value = 215.1752
$ws.Range("N3").Value = 215.1752
$ws.Range("B4").Value = 6.158
$ws.Range("C4").Value = 17.154
$ws.Range("D4").Value = 38.034
$ws.Range("E4").Value = 7.722
$ws.Range("F4").Value = 2919.128
$ws.Range("G4").Value = 2591.081
$ws.Range("H4").Value = 328.047
$ws.Range("I4").Value = 328.047
$ws.Range("K4").Value = 2731.4192
$ws.Range("L4").Value = 2591.081
$ws.Range("M4").Value = 140.3382
$ws.Range("N4").Value = 140.3382

$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = 5.84735
$ws.Range("C2").Value = 16.3438
$ws.Range("D2").Value = 35.74990000000001
$ws.Range("E2").Value = 6.2636
$ws.Range("F2").Value = 2724.69525
$ws.Range("G2").Value = 2495.7254
$ws.Range("H2").Value = 228.9699
$ws.Range("I2").Value = 228.9699
$ws.Range("K2").Value = 2650.173
$ws.Range("L2").Value = 2495.733
$ws.Range("M2").Value = 154.44
$ws.Range("N2").Value = 154.44
$ws.Range("B3").Value = 5.85
$ws.Range("C3").Value = 16.638
$ws.Range("D3").Value = 35.668
$ws.Range("E3").Value = 5.511
$ws.Range("F3").Value = 2742.579
$ws.Range("G3").Value = 2501.19
$ws.Range("H3").Value = 241.389
$ws.Range("I3").Value = 241.389
$ws.Range("K3").Value = 2651.336
$ws.Range("L3").Value = 2501.135
$ws.Range("M3").Value = 150.201
$ws.Range("N3").Value = 150.201
$ws.Range("B4").Value = 6.158
$ws.Range("C4").Value = 17.154
$ws.Range("D4").Value = 38.034
$ws.Range("E4").Value = 7.722
$ws.Range("F4").Value = 2919.128
$ws.Range("G4").Value = 2591.081
$ws.Range("H4").Value = 328.047
$ws.Range("I4").Value = 328.047
$ws.Range("K4").Value = 2683.991
$ws.Range("L4").Value = 2591.081
$ws.Range("M4").Value = 92.91
$ws.Range("N4").Value = 92.91

$ws = $wb.Worksheets.Item(5)
$ws.Range("B2").Value = 5.83895
$ws.Range("C2").Value = 15.90655
$ws.Range("D2").Value = 37.70010000000001
$ws.Range("E2").Value = 6.70675
$ws.Range("F2").Value = 2706.90085
$ws.Range("G2").Value = 2495.83625
$ws.Range("H2").Value = 211.06465
$ws.Range("I2").Value = 211.06465
$ws.Range("K2").Value = 2736.444
$ws.Range("L2").Value = 2495.841
$ws.Range("M2").Value = 240.603
$ws.Range("N2").Value = 240.603
$ws.Range("B3").Value = 5.877
$ws.Range("C3").Value = 16.197
$ws.Range("D3").Value = 36.961
$ws.Range("E3").Value = 5.974
$ws.Range("F3").Value = 2721.703
$ws.Range("G3").Value = 2500.363
$ws.Range("H3").Value = 221.34
$ws.Range("I3").Value = 221.34
$ws.Range("K3").Value = 2736.6242
$ws.Range("L3").Value = 2500.415
$ws.Range("M3").Value = 236.2092
$ws.Range("N3").Value = 236.2092
$ws.Range("B4").Value = 6.158
$ws.Range("C4").Value = 17.154
$ws.Range("D4").Value = 38.034
$ws.Range("E4").Value = 7.722
$ws.Range("F4").Value = 2919.128
$ws.Range("G4").Value = 2591.081
$ws.Range("H4").Value = 328.047
$ws.Range("I4").Value = 328.047
$ws.Range("K4").Value = 2749.1316
$ws.Range("L4").Value = 2591.081
$ws.Range("M4").Value = 158.0506
$ws.Range("N4").Value = 158.0506
